$wb = $excel.ActiveWorkbook

# --- Proximity sheet: append new row 10 ---
$wsProximity = $wb.Worksheets.Item("Proximity")
$rowProximity = $wsProximity.Range("A10:F10")
# Force text so the date-like "2026-02-01" string isn't auto-converted
# into a date serial number by the smart-paste heuristics.
$rowProximity.NumberFormat = "@"
$wsProximity.Range("A10").Value = "2026-02-01"
$wsProximity.Range("B10").Value = "13:32:32"
$wsProximity.Range("C10").Value = "13:00"
$wsProximity.Range("D10").Value = "Living Room Main Door"
$wsProximity.Range("E10").Value = "ENTER"
$wsProximity.Range("F10").Value = "User ENTERED Living Room Main Door"
# Restore the default "Normal" style so the new row matches the rest of
# the sheet (no custom formatting was applied in the source edit).
$rowProximity.Style = "Normal"

# --- Camera sheet: append new row 2 ---
$wsCamera = $wb.Worksheets.Item("Camera")
$rowCamera = $wsCamera.Range("A2:F2")
$rowCamera.NumberFormat = "@"
$wsCamera.Range("A2").Value = "2026-02-01"
$wsCamera.Range("B2").Value = "13:32:33"
$wsCamera.Range("C2").Value = "13:00"
$wsCamera.Range("D2").Value = "Living Room Main Door"
$wsCamera.Range("E2").Value = "Image Captured"
$wsCamera.Range("F2").Value = "Active"
$rowCamera.Style = "Normal"
